$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 130
$ws.Range("I12").Value = 134.28572
$ws.Range("K12").Value = 134.28572
$ws.Range("M12").Value = 35.71428
$ws.Range("H28").Value = 1511.1
$ws.Range("J28").Value = 235.33333
$ws.Range("L28").Value = 235.33333
$ws.Range("N28").Value = -1205.33333
$ws.Range("H53").Value = 921.2
$ws.Range("J53").Value = 68.333336
$ws.Range("L53").Value = 68.333336
$ws.Range("N53").Value = -1342.333336
$ws.Range("H96").Value = 20835018
$ws.Range("I96").Value = 35716804
$ws.Range("J96").Value = 520
$ws.Range("K96").Value = 107150412
$ws.Range("L96").Value = 1560
$ws.Range("M96").Value = -107149039
$ws.Range("N96").Value = -4306
$ws.Range("H99").Value = 1367.0834
$ws.Range("I99").Value = 361.6
$ws.Range("J99").Value = 2085.2856
$ws.Range("K99").Value = 1084.8
$ws.Range("L99").Value = 6255.8568
$ws.Range("M99").Value = 413.1999999999998
$ws.Range("N99").Value = -9251.856800000001
$ws.Range("H116").Value = 2387.2
$ws.Range("I116").Value = 2291.85
$ws.Range("K116").Value = 2291.85
$ws.Range("M116").Value = 1150.15
$ws.Range("H129").Value = 861.3721
$ws.Range("J129").Value = 952.6486
$ws.Range("L129").Value = 2857.9458
$ws.Range("N129").Value = -12857.9458
$ws.Range("H132").Value = 7098968.5
$ws.Range("I132").Value = 9528360
$ws.Range("J132").Value = 13243.583
$ws.Range("K132").Value = 28585080
$ws.Range("L132").Value = 39730.749
$ws.Range("M132").Value = -28582550
$ws.Range("N132").Value = -44790.749
$ws.Range("H135").Value = 1812.0952
$ws.Range("I135").Value = 473.64706
$ws.Range("J135").Value = 7500.5
$ws.Range("K135").Value = 4262.82354
$ws.Range("L135").Value = 67504.5
$ws.Range("M135").Value = -1727.82354
$ws.Range("N135").Value = -72574.5
$ws.Range("H138").Value = 1556.3541
$ws.Range("I138").Value = 1022
$ws.Range("J138").Value = 1611.6322
$ws.Range("K138").Value = 3066
$ws.Range("L138").Value = 4834.8966
$ws.Range("M138").Value = 2074
$ws.Range("N138").Value = -15114.8966
$ws.Range("H141").Value = 423.6111
$ws.Range("I141").Value = 423.6111
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1270.8333
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3909.1667
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3790.985
$ws.Range("I32").Value = 3357.0952
$ws.Range("J32").Value = 10624.75
$ws.Range("K32").Value = 3357.0952
$ws.Range("L32").Value = 10624.75
$ws.Range("M32").Value = -3070.0952
$ws.Range("N32").Value = -11198.75
$ws.Range("H61").Value = 1719.96
$ws.Range("I61").Value = 1834.4117
$ws.Range("J61").Value = 1476.75
$ws.Range("K61").Value = 1834.4117
$ws.Range("L61").Value = 1476.75
$ws.Range("M61").Value = -1622.4117
$ws.Range("N61").Value = -1900.75
$ws.Range("H74").Value = 1676.5238
$ws.Range("I74").Value = 965.9231
$ws.Range("K74").Value = 965.9231
$ws.Range("M74").Value = -91.92309999999998
$ws.Range("H77").Value = 1676.5238
$ws.Range("I77").Value = 965.9231
$ws.Range("K77").Value = 4829.6155
$ws.Range("M77").Value = -461.6154999999999
$ws.Range("H132").Value = 1191.3226
$ws.Range("I132").Value = 918.4400000000001
$ws.Range("J132").Value = 2328.3333
$ws.Range("K132").Value = 2755.32
$ws.Range("L132").Value = 6984.999899999999
$ws.Range("M132").Value = -225.3200000000002
$ws.Range("N132").Value = -12044.9999
$ws.Range("H136").Value = 1719.96
$ws.Range("I136").Value = 1834.4117
$ws.Range("J136").Value = 1476.75
$ws.Range("K136").Value = 5503.2351
$ws.Range("L136").Value = 4430.25
$ws.Range("M136").Value = -2953.2351
$ws.Range("N136").Value = -9530.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164
$ws.Range("H105").Value = 100002136
$ws.Range("I105").Value = 125002170
$ws.Range("K105").Value = 125002170
$ws.Range("M105").Value = -125000423
$ws.Range("H134").Value = 5071.5625
$ws.Range("I134").Value = 1079.6
$ws.Range("K134").Value = 3238.8
$ws.Range("M134").Value = -703.7999999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 100422.86
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 117026.664
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 117026.664
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -117726.664
$ws.Range("H43").Value = 5227.6665
$ws.Range("J43").Value = 5227.6665
$ws.Range("L43").Value = 5227.6665
$ws.Range("N43").Value = -5595.6665
$ws.Range("H99").Value = 2378.9
$ws.Range("I99").Value = 2411.125
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 2411.125
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -913.125
$ws.Range("N99").Value = -5246
$ws.Range("H101").Value = 5227.6665
$ws.Range("J101").Value = 5227.6665
$ws.Range("L101").Value = 5227.6665
$ws.Range("N101").Value = -11717.6665
$ws.Range("H107").Value = 1223.1111
$ws.Range("I107").Value = 1448.9166
$ws.Range("J107").Value = 771.5
$ws.Range("K107").Value = 1448.9166
$ws.Range("L107").Value = 771.5
$ws.Range("M107").Value = 471.0834
$ws.Range("N107").Value = -4611.5
$ws.Range("H126").Value = 2378.9
$ws.Range("I126").Value = 2411.125
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 7233.375
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -4763.375
$ws.Range("N126").Value = -11690
$ws.Range("H132").Value = 1664.9286
$ws.Range("I132").Value = 1296.7354
$ws.Range("J132").Value = 3229.75
$ws.Range("K132").Value = 3890.2062
$ws.Range("L132").Value = 9689.25
$ws.Range("M132").Value = -1360.2062
$ws.Range("N132").Value = -14749.25
$ws.Range("H134").Value = 604.9792
$ws.Range("I134").Value = 575.75
$ws.Range("K134").Value = 1727.25
$ws.Range("M134").Value = 807.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1268.4706
$ws.Range("I5").Value = 1478.9166
$ws.Range("K5").Value = 4436.7498
$ws.Range("M5").Value = -4324.7498
$ws.Range("H60").Value = 1908.9231
$ws.Range("I60").Value = 603.3333
$ws.Range("J60").Value = 2300.6
$ws.Range("K60").Value = 1809.9999
$ws.Range("L60").Value = 6901.799999999999
$ws.Range("M60").Value = -1558.9999
$ws.Range("N60").Value = -7403.799999999999
$ws.Range("H113").Value = 684.8611
$ws.Range("J113").Value = 689.8823
$ws.Range("L113").Value = 2069.6469
$ws.Range("N113").Value = -6409.6469
$ws.Range("H122").Value = 856.5625
$ws.Range("I122").Value = 491.25
$ws.Range("J122").Value = 1221.875
$ws.Range("K122").Value = 4421.25
$ws.Range("L122").Value = 10996.875
$ws.Range("M122").Value = -1971.25
$ws.Range("N122").Value = -15896.875
$ws.Range("H131").Value = 25642326
$ws.Range("I131").Value = 90909330
$ws.Range("J131").Value = 1717.2858
$ws.Range("K131").Value = 272727990
$ws.Range("L131").Value = 5151.857400000001
$ws.Range("M131").Value = -272722950
$ws.Range("N131").Value = -15231.8574
$ws.Range("H135").Value = 1268.4706
$ws.Range("I135").Value = 1478.9166
$ws.Range("K135").Value = 13310.2494
$ws.Range("M135").Value = -10775.2494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 90003360
$ws.Range("I70").Value = 83337064
$ws.Range("K70").Value = 83337064
$ws.Range("M70").Value = -83336794
$ws.Range("H73").Value = 90003360
$ws.Range("I73").Value = 83337064
$ws.Range("K73").Value = 83337064
$ws.Range("M73").Value = -83336128
$ws.Range("H102").Value = 7071.3335
$ws.Range("I102").Value = 5790
$ws.Range("K102").Value = 5790
$ws.Range("M102").Value = -4168
$ws.Range("H113").Value = 1426.5
$ws.Range("I113").Value = 1380.2858
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1380.2858
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 789.7141999999999
$ws.Range("N113").Value = -6090
$ws.Range("H126").Value = 1900
$ws.Range("I126").Value = 1680
$ws.Range("K126").Value = 5040
$ws.Range("M126").Value = -2570
$ws.Range("H132").Value = 2951.8333
$ws.Range("I132").Value = 1900
$ws.Range("K132").Value = 5700
$ws.Range("M132").Value = -3170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 875.125
$ws.Range("I16").Value = 782.5
$ws.Range("J16").Value = 1153
$ws.Range("K16").Value = 782.5
$ws.Range("L16").Value = 1153
$ws.Range("M16").Value = -612.5
$ws.Range("N16").Value = -1493
$ws.Range("H40").Value = 5175
$ws.Range("I40").Value = 3194
$ws.Range("J40").Value = 7651.25
$ws.Range("K40").Value = 3194
$ws.Range("L40").Value = 7651.25
$ws.Range("M40").Value = -3058
$ws.Range("N40").Value = -7923.25
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960
$ws.Range("H136").Value = 1328.4445
$ws.Range("I136").Value = 1122.3334
$ws.Range("J136").Value = 2049.8333
$ws.Range("K136").Value = 3367.0002
$ws.Range("L136").Value = 6149.499899999999
$ws.Range("M136").Value = -817.0001999999999
$ws.Range("N136").Value = -11249.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 38255.5
$ws.Range("J20").Value = 38255.5
$ws.Range("L20").Value = 38255.5
$ws.Range("N20").Value = -38735.5
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = ""
$ws.Range("H107").Value = 391.30435
$ws.Range("J107").Value = 429.8
$ws.Range("L107").Value = 1289.4
$ws.Range("M107").Value = 429.8
$ws.Range("N107").Value = -5129.4
